$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="70.967.10"'
$ws.Range("E2").Formula = '="  +6.03%  "'
$ws.Range("D3").Formula = '="3.667.32"'
$ws.Range("E3").Formula = '="  +18.15%  "'
$ws.Range("E4").Formula = '="  -0.01%  "'
$ws.Range("D5").Formula = '="619.66"'
$ws.Range("E5").Formula = '="  +7.50%  "'
$ws.Range("D6").Formula = '="182.06"'
$ws.Range("E6").Formula = '="  +2.55%  "'
$ws.Range("D7").Formula = '="3.662.36"'
$ws.Range("E7").Formula = '="  +18.00%  "'
$ws.Range("E8").Formula = '="  -0.04%  "'
$ws.Range("D9").Formula = '="0.542"'
$ws.Range("E9").Formula = '="  +5.84%  "'
$ws.Range("E10").Formula = '="  +7.97%  "'
$ws.Range("D11").Formula = '="6.69"'
$ws.Range("E11").Formula = '="  +5.46%  "'
$ws.Range("D12").Formula = '="0.501"'
$ws.Range("E12").Formula = '="  +7.22%  "'
$ws.Range("D13").Formula = '="40.35"'
$ws.Range("E13").Formula = '="  +11.66%  "'
$ws.Range("D14").Formula = '="0.0000255"'
$ws.Range("E14").Formula = '="  +6.43%  "'
$ws.Range("D15").Formula = '="4.277.94"'
$ws.Range("E15").Formula = '="  +18.09%  "'
$ws.Range("D16").Formula = '="70.998.02"'
$ws.Range("D17").Formula = '="3.653.86"'
$ws.Range("E17").Formula = '="  +17.70%  "'
$ws.Range("E18").Formula = '="  +2.10%  "'
$ws.Range("E19").Formula = '="  +7.24%  "'
$ws.Range("D20").Formula = '="519.45"'
$ws.Range("E20").Formula = '="  +8.32%  "'
$ws.Range("D21").Formula = '="16.90"'
$ws.Range("E21").Formula = '="  +1.03%  "'
$ws.Range("D22").Formula = '="9.22"'
$ws.Range("E22").Formula = '="  +18.62%  "'
$ws.Range("E23").Formula = '="  +7.62%  "'
$ws.Range("D24").Formula = '="2.54"'
$ws.Range("E24").Formula = '="  +12.85%  "'
$ws.Range("D25").Formula = '="88.55"'
$ws.Range("E25").Formula = '="  +5.93%  "'
$ws.Range("D26").Formula = '="13.51"'
$ws.Range("E26").Formula = '="  +7.49%  "'
$ws.Range("D27").Formula = '="11.07"'
$ws.Range("E27").Formula = '="  +9.43%  "'
$ws.Range("E28").Formula = '="  -0.02%  "'
$ws.Range("E29").Formula = '="  +10.94%  "'
$ws.Range("D30").Formula = '="8.17"'
$ws.Range("E30").Formula = '="  +3.31%  "'
$ws.Range("D31").Formula = '="2.91"'
$ws.Range("E31").Formula = '="  +11.83%  "'
$ws.Range("E32").Formula = '="  +18.18%  "'
$ws.Range("E33").Formula = '="  +12.90%  "'
$ws.Range("E34").Formula = '="  +4.59%  "'
$ws.Range("E35").Formula = '="  -0.08%  "'
$ws.Range("E36").Formula = '="  +9.58%  "'
$ws.Range("E37").Formula = '="  +9.74%  "'
$ws.Range("D38").Formula = '="0.348"'
$ws.Range("E38").Formula = '="  +11.60%  "'
$ws.Range("E39").Formula = '="  +9.66%  "'
$ws.Range("E40").Formula = '="  +6.96%  "'
$ws.Range("D41").Formula = '="51.33"'
$ws.Range("E41").Formula = '="  +4.63%  "'
$ws.Range("D42").Formula = '="45.58"'
$ws.Range("E42").Formula = '="  -5.67%  "'
$ws.Range("D43").Formula = '="432.78"'
$ws.Range("E43").Formula = '="  +16.16%  "'
$ws.Range("E44").Formula = '="  +6.02%  "'
$ws.Range("D45").Formula = '="3.110.36"'
$ws.Range("E45").Formula = '="  +11.16%  "'
$ws.Range("D46").Formula = '="2.85"'
$ws.Range("E46").Formula = '="  +4.94%  "'
$ws.Range("D47").Formula = '="0.0369"'
$ws.Range("E47").Formula = '="  +7.29%  "'
$ws.Range("D48").Formula = '="28.33"'
$ws.Range("E48").Formula = '="  +10.01%  "'
$ws.Range("D49").Formula = '="140.02"'
$ws.Range("E49").Formula = '="  +3.11%  "'
$ws.Range("D51").Formula = '="2.48"'
$ws.Range("E51").Formula = '="  +10.89%  "'

$src = $ws.Range("D2:E51")
$src.Copy()
$src.PasteSpecial(-4163)
$excel.CutCopyMode = 0

